# Update gh-pages to output generated at 456a3b4
# This script updates the "F" column (想去人数 / want-to-go count) values
# on three of the four worksheets to reflect newly generated data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 347
$ws1.Range("F7").Value  = 5982
$ws1.Range("F8").Value  = 10193
$ws1.Range("F9").Value  = 4030
$ws1.Range("F15").Value = 26
$ws1.Range("F16").Value = 153
$ws1.Range("F18").Value = 5675
$ws1.Range("F23").Value = 8515
$ws1.Range("F30").Value = 199
$ws1.Range("F31").Value = 1877
$ws1.Range("F36").Value = 275
$ws1.Range("F42").Value = 78
$ws1.Range("F44").Value = 1415
$ws1.Range("F45").Value = 2285

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F19").Value = 3

# --- Sheet "本地生活" ---
# (no changes in this sheet)

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 347
$ws4.Range("F7").Value  = 5982
$ws4.Range("F8").Value  = 4030
$ws4.Range("F16").Value = 26
$ws4.Range("F17").Value = 153
$ws4.Range("F20").Value = 5675
$ws4.Range("F25").Value = 8515
$ws4.Range("F31").Value = 199
$ws4.Range("F32").Value = 1877
$ws4.Range("F36").Value = 275
$ws4.Range("F42").Value = 78
$ws4.Range("F44").Value = 1415
$ws4.Range("F46").Value = 2285
